$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reproduce the worksheet's stale outline-level watermarks -------------
# The original workbook already carried outlineLevelRow="1"/outlineLevelCol="1"
# in sheetFormatPr despite no row/col currently showing an outlineLevel
# (a leftover from earlier grouping that was later removed). The target
# revision bumps these to outlineLevelRow="3" / outlineLevelCol="2" while
# still leaving no visible grouped row/col in the final sheet, so we
# reproduce that by grouping a throwaway column/row to the desired depth
# and then deleting it again.
$ws.Columns("D:D").Group()
$ws.Columns("D:D").Group()
$ws.Columns("D:D").Delete()

$ws.Rows("9:9").Group()
$ws.Rows("9:9").Group()
$ws.Rows("9:9").Group()
$ws.Rows("9:9").Delete()

# --- Expand the sample table ------------------------------------------------
# Before:
#   A1=镜头   B1=资产
#   A2=SDKTEST_EP01_01_sc001  B2=asset1
# After:
#   A1=镜头   B1=资产1        C1=资产2
#   A2=SDKTEST_EP01_01_sc001  B2=asset1   C2=asset2
#   A3=SDKTEST_EP01_01_sc002  B3=asset1
#   A4=SDKTEST_EP01_01_sc003              C4=asset2
$ws.Range("B1").Value = "资产1"
$ws.Range("C1").Value = "资产2"

$ws.Range("B2").Value = "asset1"
$ws.Range("C2").Value = "asset2"

$ws.Range("A3").Value = "SDKTEST_EP01_01_sc002"
$ws.Range("B3").Value = "asset1"

$ws.Range("A4").Value = "SDKTEST_EP01_01_sc003"
$ws.Range("C4").Value = "asset2"

# --- Match the saved cursor position ---------------------------------------
$ws.Range("A9").Select()
